# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 09:55"

# Row 62 <-> 63 : Singapur / Armenia swap places (Armenia overtakes Singapur) with refreshed data
$ws.Cells.Item(62, 1).Value = "Armenia"
$ws.Cells.Item(62, 2).Value = 58624
$ws.Cells.Item(62, 3).Value = 1058
$ws.Cells.Item(62, 4).Value = 46713
$ws.Cells.Item(62, 5).Value = 10872
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 7
$ws.Cells.Item(62, 8).Value = 1039

$ws.Cells.Item(63, 1).Value = "Singapur"
$ws.Cells.Item(63, 2).Value = 57889
$ws.Cells.Item(63, 3).Value = 5
$ws.Cells.Item(63, 4).Value = 57740
$ws.Cells.Item(63, 5).Value = 121
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 28

# Row 75 <-> 76 : Afganistan / Hungria swap places (Hungria overtakes Afganistan) with refreshed data
$ws.Cells.Item(75, 1).Value = "Hungria"
$ws.Cells.Item(75, 2).Value = 40782
$ws.Cells.Item(75, 3).Value = 920
$ws.Cells.Item(75, 4).Value = 12164
$ws.Cells.Item(75, 5).Value = 27595
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 27
$ws.Cells.Item(75, 8).Value = 1023

$ws.Cells.Item(76, 1).Value = "Afganistan"
$ws.Cells.Item(76, 2).Value = 39994
$ws.Cells.Item(76, 3).Value = 66
$ws.Cells.Item(76, 4).Value = 33354
$ws.Cells.Item(76, 5).Value = 5160
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 1480

# Row 101 : Georgia, refreshed data (no name swap)
$ws.Cells.Item(101, 1).Value = "Georgia"
$ws.Cells.Item(101, 2).Value = 13521
$ws.Cells.Item(101, 3).Value = 680
$ws.Cells.Item(101, 4).Value = 7159
$ws.Cells.Item(101, 5).Value = 6253
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 7
$ws.Cells.Item(101, 8).Value = 109

# Row 140 <-> 141 : Malta / Estonia swap places (Estonia overtakes Malta) with refreshed data
$ws.Cells.Item(140, 1).Value = "Estonia"
$ws.Cells.Item(140, 2).Value = 3947
$ws.Cells.Item(140, 3).Value = 39
$ws.Cells.Item(140, 4).Value = 3060
$ws.Cells.Item(140, 5).Value = 819
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 68

$ws.Cells.Item(141, 1).Value = "Malta"
$ws.Cells.Item(141, 2).Value = 3937
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 3012
$ws.Cells.Item(141, 5).Value = 881
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 44

# Row 194 : Bonaire, San Eustaquio y Saba, refreshed data (no name swap)
$ws.Cells.Item(194, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(194, 2).Value = 150
$ws.Cells.Item(194, 3).Value = 2
$ws.Cells.Item(194, 4).Value = 111
$ws.Cells.Item(194, 5).Value = 37
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 2
